$p = $ppt.ActivePresentation

# Slide 6: title was "Crisper Knockout (Single Perturbation)" + line break + "Prediction Horizon=1"
#          -> becomes just "Crisper Knockout" (single run, no break, no second line)
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(1).TextFrame.TextRange
$tr6.Characters(1, $tr6.Length).Text = "Crisper Knockout"

# Slide 18: title was "Drug Development(Single Perturbation)" + line break + "Prediction Horizon=1"
#           -> becomes just "Drug Development" (single run, no break, no second line)
$s18 = $p.Slides.Item(18)
$tr18 = $s18.Shapes.Item(1).TextFrame.TextRange
$tr18.Characters(1, $tr18.Length).Text = "Drug Development"
